$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 425.625
$ws.Range("I33").Value = 144.54546
$ws.Range("J33").Value = 1044
$ws.Range("K33").Value = 144.54546
$ws.Range("L33").Value = 1044
$ws.Range("M33").Value = 84.45454000000001
$ws.Range("N33").Value = -1502

$ws.Range("H137").Value = 23817.562
$ws.Range("I137").Value = 24099
$ws.Range("J137").Value = 21397.2
$ws.Range("K137").Value = 72297
$ws.Range("L137").Value = 64191.60000000001
$ws.Range("M137").Value = -69747
$ws.Range("N137").Value = -69291.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1336
$ws.Range("I2").Value = 1170.3077
$ws.Range("J2").Value = 1695
$ws.Range("K2").Value = 1170.3077
$ws.Range("L2").Value = 1695
$ws.Range("M2").Value = -1057.3077
$ws.Range("N2").Value = -1921

$ws.Range("H61").Value = 1716.3529
$ws.Range("I61").Value = 598.5333000000001
$ws.Range("J61").Value = 10100
$ws.Range("K61").Value = 598.5333000000001
$ws.Range("L61").Value = 10100
$ws.Range("M61").Value = -386.5333000000001
$ws.Range("N61").Value = -10524

$ws.Range("H74").Value = 50432.293
$ws.Range("I74").Value = 72940.21000000001
$ws.Range("J74").Value = 1953.6923
$ws.Range("K74").Value = 72940.21000000001
$ws.Range("L74").Value = 1953.6923
$ws.Range("M74").Value = -72066.21000000001
$ws.Range("N74").Value = -3701.6923

$ws.Range("H77").Value = 50432.293
$ws.Range("I77").Value = 72940.21000000001
$ws.Range("J77").Value = 1953.6923
$ws.Range("K77").Value = 364701.05
$ws.Range("L77").Value = 9768.461499999999
$ws.Range("M77").Value = -360333.05
$ws.Range("N77").Value = -18504.4615

$ws.Range("H116").Value = 1336
$ws.Range("I116").Value = 1170.3077
$ws.Range("J116").Value = 1695
$ws.Range("K116").Value = 1170.3077
$ws.Range("L116").Value = 1695
$ws.Range("M116").Value = 1123.6923
$ws.Range("N116").Value = -6283

$ws.Range("H132").Value = 2337316
$ws.Range("I132").Value = 2685651.8
$ws.Range("J132").Value = 1013640.2
$ws.Range("K132").Value = 8056955.399999999
$ws.Range("L132").Value = 3040920.6
$ws.Range("M132").Value = -8054425.399999999
$ws.Range("N132").Value = -3045980.6

$ws.Range("H136").Value = 1716.3529
$ws.Range("I136").Value = 598.5333000000001
$ws.Range("J136").Value = 10100
$ws.Range("K136").Value = 1795.5999
$ws.Range("L136").Value = 30300
$ws.Range("M136").Value = 754.4000999999998
$ws.Range("N136").Value = -35400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1336
$ws.Range("I3").Value = 1170.3077
$ws.Range("J3").Value = 1695
$ws.Range("K3").Value = 1170.3077
$ws.Range("L3").Value = 1695
$ws.Range("M3").Value = -1056.3077
$ws.Range("N3").Value = -1923

$ws.Range("H68").Value = 20000
$ws.Range("J68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("N68").Value = -21622

$ws.Range("H71").Value = 20000
$ws.Range("J71").Value = 20000
$ws.Range("L71").Value = 60000
$ws.Range("N71").Value = -68112

$ws.Range("H134").Value = 50833
$ws.Range("I134").Value = 1939.5834
$ws.Range("J134").Value = 116024.22
$ws.Range("K134").Value = 5818.7502
$ws.Range("L134").Value = 348072.66
$ws.Range("M134").Value = -3283.7502
$ws.Range("N134").Value = -353142.66

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H31").Value = 28137.44
$ws.Range("I31").Value = 35082.086
$ws.Range("K31").Value = 35082.086
$ws.Range("M31").Value = -34787.086

$ws.Range("H34").Value = 28137.44
$ws.Range("I34").Value = 35082.086
$ws.Range("K34").Value = 35082.086
$ws.Range("M34").Value = -34880.086

$ws.Range("H53").Value = 19500
$ws.Range("J53").Value = 19500
$ws.Range("L53").Value = 19500
$ws.Range("N53").Value = -20714

$ws.Range("H58").Value = 6577.353
$ws.Range("I58").Value = 755
$ws.Range("J58").Value = 25500
$ws.Range("K58").Value = 755
$ws.Range("L58").Value = 25500
$ws.Range("M58").Value = -552
$ws.Range("N58").Value = -25906

$ws.Range("H132").Value = 1812.4166
$ws.Range("I132").Value = 1293.9474
$ws.Range("J132").Value = 3782.6
$ws.Range("K132").Value = 3881.8422
$ws.Range("L132").Value = 11347.8
$ws.Range("M132").Value = -1351.8422
$ws.Range("N132").Value = -16407.8

$ws.Range("H134").Value = 10418021
$ws.Range("I134").Value = 1201.75
$ws.Range("J134").Value = 41668476
$ws.Range("K134").Value = 3605.25
$ws.Range("L134").Value = 125005428
$ws.Range("M134").Value = -1070.25
$ws.Range("N134").Value = -125010498

$ws.Range("H136").Value = 6577.353
$ws.Range("I136").Value = 755
$ws.Range("J136").Value = 25500
$ws.Range("K136").Value = 2265
$ws.Range("L136").Value = 76500
$ws.Range("M136").Value = 285
$ws.Range("N136").Value = -81600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 10871.444
$ws.Range("I110").Value = 1950
$ws.Range("J110").Value = 11986.625
$ws.Range("K110").Value = 5850
$ws.Range("L110").Value = 35959.875
$ws.Range("M110").Value = -1760
$ws.Range("N110").Value = -44139.875

$ws.Range("H114").Value = 1217.5555
$ws.Range("I114").Value = 614
$ws.Range("J114").Value = 1390
$ws.Range("K114").Value = 1842
$ws.Range("L114").Value = 4170
$ws.Range("M114").Value = 1412
$ws.Range("N114").Value = -10678

$ws.Range("H115").Value = 2284
$ws.Range("I115").Value = 775
$ws.Range("J115").Value = 2661.25
$ws.Range("K115").Value = 2325
$ws.Range("L115").Value = 7983.75
$ws.Range("M115").Value = -1150
$ws.Range("N115").Value = -10333.75

$ws.Range("H125").Value = 7163.3335
$ws.Range("J125").Value = 9250
$ws.Range("L125").Value = 27750
$ws.Range("N125").Value = -37590

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 47005.332
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 70008
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 70008
$ws.Range("M9").Value = -830
$ws.Range("N9").Value = -70348

$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4996

$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24984

$ws.Range("H132").Value = 34225.16
$ws.Range("I132").Value = 1521.3077
$ws.Range("J132").Value = 204285.2
$ws.Range("K132").Value = 4563.9231
$ws.Range("L132").Value = 612855.6000000001
$ws.Range("M132").Value = -2033.9231
$ws.Range("N132").Value = -617915.6000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 1000207
$ws.Range("J69").Value = 1000207
$ws.Range("L69").Value = 1000207
$ws.Range("N69").Value = -1001829

$ws.Range("H72").Value = 1000207
$ws.Range("J72").Value = 1000207
$ws.Range("L72").Value = 3000621
$ws.Range("N72").Value = -3008733

$ws.Range("H132").Value = 390515.4
$ws.Range("I132").Value = 114186.945
$ws.Range("J132").Value = 773124.0600000001
$ws.Range("K132").Value = 342560.835
$ws.Range("L132").Value = 2319372.18
$ws.Range("M132").Value = -340030.835
$ws.Range("N132").Value = -2324432.18

$ws.Range("H136").Value = 401486.8
$ws.Range("I136").Value = 588985.3
$ws.Range("J136").Value = 3052.5
$ws.Range("K136").Value = 1766955.9
$ws.Range("L136").Value = 9157.5
$ws.Range("M136").Value = -1764405.9
$ws.Range("N136").Value = -14257.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7039.421
$ws.Range("I132").Value = 1349.7333
$ws.Range("K132").Value = 4049.199900000001
$ws.Range("M132").Value = -1519.199900000001

$ws.Range("H136").Value = 3017758.8
$ws.Range("I136").Value = 3761602.5
$ws.Range("J136").Value = 1251130
$ws.Range("K136").Value = 11284807.5
$ws.Range("L136").Value = 3753390
$ws.Range("M136").Value = -11282257.5
$ws.Range("N136").Value = -3758490
